$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "afleveres uge 16 2012!" -> "afleveres "
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("afleveres uge 16 2012!") | Out-Null
$r1.Text = "afleveres "

# ---------------------------------------------------------------------------
# 2) "(ændres)" -> "fredag d. 20 kl. 9:00 2012"
#    (this run sits right after the one we just edited)
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("(ændres)") | Out-Null
$r2.Text = "fredag d. 20 kl. 9:00 2012"

# The two edits above collapse into a single run because the surrounding
# runs share identical formatting. Re-split the new date text back into its
# own run (matching the original authoring) by nudging a formatting
# property off and back onto its original value.
$r3 = $d.Content
$r3.Find.Execute("fredag d. 20 kl. 9:00 2012") | Out-Null
$origColor = $r3.Font.Color
$r3.Font.Color = wdColorWhite
$r3.Font.Color = $origColor

# ---------------------------------------------------------------------------
# 3) Insert a blank paragraph before the page-break paragraph, carrying the
#    same line-spacing used by the paragraph above it, and mark the
#    following page break with a lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$brPara = $d.Paragraphs(8)
$brRange = $brPara.Range
Write-Output ("brPara text=[" + $brRange.Text + "]")

$insertPoint = $d.Range($brPara.Range.Start, $brPara.Range.Start)
$insertPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs(8)
$newPara.Format.LineSpacingRule = 5
$newPara.Format.LineSpacing = $word.LinesToPoints(1.5)

$word.Selection.GoTo(3, 0, 0, 10) | Out-Null
